$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.033697341275478
$ws.Range("D2").Value = 1.035959640304085
$ws.Range("E2").Value = 1.037340588814893
$ws.Range("F2").Value = 1.032337102294503
$ws.Range("I2").Value = 1.035876791816523
$ws.Range("J2").Value = 1.038820645446166
$ws.Range("K2").Value = 1.038754751167188
$ws.Range("L2").Value = 1.040131751611263
$ws.Range("M2").Value = 1.035142624361614
$ws.Range("N2").Value = 1.016827439257092
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.034837870960334
$ws.Range("D3").Value = 1.036800500449208
$ws.Range("E3").Value = 1.038429599452509
$ws.Range("F3").Value = 1.034109656327795
$ws.Range("I3").Value = 1.036180491228094
$ws.Range("J3").Value = 1.039603029484229
$ws.Range("K3").Value = 1.039405030712471
$ws.Range("L3").Value = 1.041029810551273
$ws.Range("M3").Value = 1.036721352937739
$ws.Range("N3").Value = 1.017093817305043
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.035575160368156
$ws.Range("D4").Value = 1.037343983904056
$ws.Range("E4").Value = 1.039133911551767
$ws.Range("F4").Value = 1.035255747469733
$ws.Range("I4").Value = 1.036375484373667
$ws.Range("J4").Value = 1.040108053687531
$ws.Range("K4").Value = 1.03982455454903
$ws.Range("L4").Value = 1.041609976657832
$ws.Range("M4").Value = 1.037741595764192
$ws.Range("N4").Value = 1.017265599945949
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.035884949868084
$ws.Range("D5").Value = 1.037572319865443
$ws.Range("E5").Value = 1.039429922722727
$ws.Range("F5").Value = 1.035737365516398
$ws.Range("I5").Value = 1.036457096293951
$ws.Range("J5").Value = 1.040320073332185
$ws.Range("K5").Value = 1.040000624514394
$ws.Range("L5").Value = 1.041853656076519
$ws.Range("M5").Value = 1.03817020304926
$ws.Range("N5").Value = 1.017337678666173
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.03593695518136
$ws.Range("D6").Value = 1.037610650045421
$ws.Range("E6").Value = 1.039479619580258
$ws.Range("F6").Value = 1.035818219933011
$ws.Range("I6").Value = 1.036470778021135
$ws.Range("J6").Value = 1.040355655269368
$ws.Range("K6").Value = 1.040030170013881
$ws.Range("L6").Value = 1.04189455794346
$ws.Range("M6").Value = 1.038242150587083
$ws.Range("N6").Value = 1.017349772883505
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.035579300443078
$ws.Range("D7").Value = 1.037347035506369
$ws.Range("E7").Value = 1.039137867185445
$ws.Range("F7").Value = 1.035262183642427
$ws.Range("I7").Value = 1.036376576302204
$ws.Range("J7").Value = 1.040110887851281
$ws.Range("K7").Value = 1.03982690837302
$ws.Range("L7").Value = 1.041613233585404
$ws.Range("M7").Value = 1.037747324015048
$ws.Range("N7").Value = 1.017266563609633
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.034082937294239
$ws.Range("D8").Value = 1.036243939919204
$ws.Range("E8").Value = 1.037708699041882
$ws.Range("F8").Value = 1.032936331109785
$ws.Range("I8").Value = 1.035979743838795
$ws.Range("J8").Value = 1.039085311848311
$ws.Range("K8").Value = 1.038974776029039
$ws.Range("L8").Value = 1.040435450377758
$ws.Range("M8").Value = 1.035676437374622
$ws.Range("N8").Value = 1.016917583771953
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031440589040341
$ws.Range("D9").Value = 1.034295415358004
$ws.Range("E9").Value = 1.035187533109893
$ws.Range("F9").Value = 1.028830790322248
$ws.Range("I9").Value = 1.03526878897738
$ws.Range("J9").Value = 1.037268597441004
$ws.Range("K9").Value = 1.037463563754977
$ws.Range("L9").Value = 1.03835276353889
$ws.Range("M9").Value = 1.032016933725559
$ws.Range("N9").Value = 1.016298154984115
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.029675096370383
$ws.Range("D10").Value = 1.032993128676413
$ws.Range("E10").Value = 1.033504732571351
$ws.Range("F10").Value = 1.026088393757275
$ws.Range("I10").Value = 1.03478690696604
$ws.Range("J10").Value = 1.036050921637712
$ws.Range("K10").Value = 1.036449504960099
$ws.Range("L10").Value = 1.036959270196253
$ws.Range("M10").Value = 1.029569783176021
$ws.Range("N10").Value = 1.015882151716184
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028909645024494
$ws.Range("D11").Value = 1.032428427436566
$ws.Range("E11").Value = 1.032775547660807
$ws.Range("F11").Value = 1.024899498116727
$ws.Range("I11").Value = 1.034576357304364
$ws.Range("J11").Value = 1.035522074719935
$ws.Range("K11").Value = 1.036008822102219
$ws.Range("L11").Value = 1.036354646102886
$ws.Range("M11").Value = 1.028508246856749
$ws.Range("N11").Value = 1.015701285452943
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028625171151468
$ws.Range("D12").Value = 1.032218549770686
$ws.Range("E12").Value = 1.032504614605495
$ws.Range("F12").Value = 1.024457664425698
$ws.Range("I12").Value = 1.034497864433361
$ws.Range("J12").Value = 1.035325396555764
$ws.Range("K12").Value = 1.035844892149048
$ws.Range("L12").Value = 1.03612987382729
$ws.Range("M12").Value = 1.028113648674378
$ws.Range("N12").Value = 1.01563399259304
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028686198663313
$ws.Range("D13").Value = 1.032263574815988
$ws.Range("E13").Value = 1.032562734393766
$ws.Range("F13").Value = 1.024552449617227
$ws.Range("I13").Value = 1.034514714352224
$ws.Range("J13").Value = 1.035367595633246
$ws.Range("K13").Value = 1.035880066606327
$ws.Range("L13").Value = 1.036178096794982
$ws.Range("M13").Value = 1.028198304899221
$ws.Range("N13").Value = 1.01564843217908
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028886133430676
$ws.Range("D14").Value = 1.032411081392975
$ws.Range("E14").Value = 1.032753153926203
$ws.Range("F14").Value = 1.024862980661574
$ws.Range("I14").Value = 1.034569874887026
$ws.Range("J14").Value = 1.03550582217097
$ws.Range("K14").Value = 1.035995276523558
$ws.Range("L14").Value = 1.03633607020738
$ws.Range("M14").Value = 1.028475635352606
$ws.Range("N14").Value = 1.015695725274657
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.029009299697391
$ws.Range("D15").Value = 1.032501948775716
$ws.Range("E15").Value = 1.032870466819526
$ws.Range("F15").Value = 1.025054278935046
$ws.Range("I15").Value = 1.034603823266136
$ws.Range("J15").Value = 1.035590956103749
$ws.Range("K15").Value = 1.036066229197323
$ws.Range("L15").Value = 1.036433377848098
$ws.Range("M15").Value = 1.028646468222967
$ws.Range("N15").Value = 1.015724849369323
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029725875861851
$ws.Range("D16").Value = 1.03303058896703
$ws.Range("E16").Value = 1.033553114935307
$ws.Range("F16").Value = 1.026167265942798
$ws.Range("I16").Value = 1.034800840501232
$ws.Range("J16").Value = 1.036085985835978
$ws.Range("K16").Value = 1.036478717984672
$ws.Range("L16").Value = 1.036999370902275
$ws.Range("M16").Value = 1.029640192855148
$ws.Range("N16").Value = 1.015894139680799
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030175099797833
$ws.Range("D17").Value = 1.033361974792362
$ws.Range("E17").Value = 1.033981180411138
$ws.Range("F17").Value = 1.026865024993354
$ws.Range("I17").Value = 1.034923916960961
$ws.Range("J17").Value = 1.036396078586595
$ws.Range("K17").Value = 1.036737034531195
$ws.Range("L17").Value = 1.037354071484195
$ws.Range("M17").Value = 1.03026301307928
$ws.Range("N17").Value = 1.016000133927908
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030437030085237
$ws.Range("D18").Value = 1.033555189238146
$ws.Range("E18").Value = 1.034230813681112
$ws.Range("F18").Value = 1.027271879837248
$ws.Range("I18").Value = 1.034995522948101
$ws.Range("J18").Value = 1.036576797671738
$ws.Range("K18").Value = 1.036887553054157
$ws.Range("L18").Value = 1.037560843628651
$ws.Range("M18").Value = 1.030626110570095
$ws.Range("N18").Value = 1.016061887802518
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030526325638482
$ws.Range("D19").Value = 1.033621057316902
$ws.Range("E19").Value = 1.034315923771008
$ws.Range("F19").Value = 1.027410584071765
$ws.Range("I19").Value = 1.035019907830952
$ws.Range("J19").Value = 1.036638392378315
$ws.Range("K19").Value = 1.036938850065067
$ws.Range("L19").Value = 1.037631327537112
$ws.Range("M19").Value = 1.030749886699376
$ws.Range("N19").Value = 1.016082932296508
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030126912116845
$ws.Range("D20").Value = 1.033326428239661
$ws.Range("E20").Value = 1.033935258239641
$ws.Range("F20").Value = 1.026790176225348
$ws.Range("I20").Value = 1.034910730900043
$ws.Range("J20").Value = 1.036362824391346
$ws.Range("K20").Value = 1.036709335466759
$ws.Range("L20").Value = 1.037316027768187
$ws.Range("M20").Value = 1.030196209369458
$ws.Range("N20").Value = 1.015988769075055
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028827261827747
$ws.Range("D21").Value = 1.032367647779137
$ws.Range("E21").Value = 1.032697082388926
$ws.Range("F21").Value = 1.024771543311103
$ws.Range("I21").Value = 1.03455363937088
$ws.Range("J21").Value = 1.035465124571894
$ws.Range("K21").Value = 1.03596135673589
$ws.Range("L21").Value = 1.036289556186453
$ws.Range("M21").Value = 1.028393976709361
$ws.Range("N21").Value = 1.015681801709372
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.028009243672713
$ws.Range("D22").Value = 1.031764114843136
$ws.Range("E22").Value = 1.031918120129035
$ws.Range("F22").Value = 1.023501041589965
$ws.Range("I22").Value = 1.034327470211495
$ws.Range("J22").Value = 1.03489931009939
$ws.Range("K22").Value = 1.035489679078344
$ws.Range("L22").Value = 1.035643083828226
$ws.Range("M22").Value = 1.027259121900553
$ws.Range("N22").Value = 1.015488156092582
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028442974900137
$ws.Range("D23").Value = 1.032084126979107
$ws.Range("E23").Value = 1.032331108510677
$ws.Range("F23").Value = 1.024174686384034
$ws.Range("I23").Value = 1.034447523700911
$ws.Range("J23").Value = 1.035199392192832
$ws.Range("K23").Value = 1.035739857123704
$ws.Range("L23").Value = 1.035985895226966
$ws.Range("M23").Value = 1.027860896356801
$ws.Range("N23").Value = 1.015590872516418
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.030148686347297
$ws.Range("D24").Value = 1.033342490435568
$ws.Range("E24").Value = 1.03395600864689
$ws.Range("F24").Value = 1.026823997581606
$ws.Range("I24").Value = 1.034916689677023
$ws.Range("J24").Value = 1.036377851002683
$ws.Range("K24").Value = 1.03672185195482
$ws.Range("L24").Value = 1.037333218452544
$ws.Range("M24").Value = 1.030226395654959
$ws.Range("N24").Value = 1.015993904581532
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.032124377919633
$ws.Range("D25").Value = 1.03479972536272
$ws.Range("E25").Value = 1.03583966061755
$ws.Range("F25").Value = 1.029893073103564
$ws.Range("I25").Value = 1.03545397831206
$ws.Range("J25").Value = 1.037739403292432
$ws.Range("K25").Value = 1.037855401546925
$ws.Range("L25").Value = 1.038892065264326
$ws.Range("M25").Value = 1.032964282632046
$ws.Range("N25").Value = 1.016458827113751
